# Adds the new CFR method rows (264-272) to Sheet1, matching the rows
# already present in the list (columns A: Char_Name, B: CASNumber,
# C: Method_Code, D: Method_Context, E: CFR_Method).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 264; A = "Total suspended solids";                             C = "2540-D" },
    @{ Row = 265; A = "Total dissolved solids";                             C = "2540-C" },
    @{ Row = 266; A = "Dissolved oxygen (DO)";                              C = "4500-O-G" },
    @{ Row = 267; A = "Biochemical oxygen demand, standard conditions";     C = "5210-B" },
    @{ Row = 268; A = "Chemical oxygen demand";                             C = "5220-C" },
    @{ Row = 269; A = "Chemical oxygen demand";                             C = "5220-D" },
    @{ Row = 270; A = "Alkalinity, bicarbonate";                            C = "2320-B" },
    @{ Row = 271; A = "Alkalinity, carbonate";                              C = "2320-B" },
    @{ Row = 272; A = "Alkalinity, Hydroxide";                              C = "2320-B" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $cellA = $ws.Cells.Item($rowNum, 1)
    $cellA.HorizontalAlignment = -4131
    $cellA.Value = $r.A

    $cellC = $ws.Cells.Item($rowNum, 3)
    $cellC.HorizontalAlignment = -4131
    $cellC.NumberFormat = "@"
    $cellC.Value = $r.C

    $cellD = $ws.Cells.Item($rowNum, 4)
    $cellD.Value = "American Public Health Association (Standards Methods - SM)"

    $cellE = $ws.Cells.Item($rowNum, 5)
    $cellE.HorizontalAlignment = -4131
    $cellE.Value = "Yes"
}

# Match the author's final selection (the last edited block, D270:E272).
$ws.Range("D270:E272").Select() | Out-Null
